$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'231.05"
$ws.Range("D3").Value = "'22.77"
$ws.Range("D4").Value = "'5.288"
$ws.Range("D5").Value = "'0.05602"
$ws.Range("D6").Value = "'3.374"
$ws.Range("D7").Value = "'6.463"
$ws.Range("D8").Value = "'1.063"
$ws.Range("D9").Value = "'0.7823"
$ws.Range("D10").Value = "'0.1387"
$ws.Range("D11").Value = "'0.07403"
$ws.Range("D13").Value = "'0.02965"
$ws.Range("D14").Value = "'0.09259"
$ws.Range("D15").Value = "'0.001642"
$ws.Range("D16").Value = "'3.264"
$ws.Range("D17").Value = "'0.04763"
$ws.Range("D18").Value = "'0.0005786"
$ws.Range("D19").Value = "'0.006246"
$ws.Range("D20").Value = "'0.005238"
$ws.Range("D21").Value = "'0.001053"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.979"
$ws.Range("D27").Value = "'0.0004996"
$ws.Range("E27").Value = "'26UpBotsUBXTBestin24h"
$ws.Range("D40").Value = "'0.04038"
$ws.Range("D41").Value = "'0.007000"
$ws.Range("E41").Value = "'40KickTokenKICK"
$ws.Range("D43").Value = "'0.003210"
$ws.Range("D44").Value = "'0.009280"
$ws.Range("D45").Value = "'0.00005436"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.7845"
$ws.Range("D48").Value = "'0.04125"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.01009"
